$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: update title and link
$ws.Range("D5").Value = "신호 처리 서론"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/01/03/signal_processing_introduction.html"

# Row 41: update title and link
$ws.Range("D41").Value = "유튜브가 딥러닝을 도입한 방식 (1)"
$ws.Range("E41").Value = "http://cloudinsight.net/ai/%ec%9c%a0%ed%8a%9c%eb%b8%8c%ea%b0%80-%eb%94%a5%eb%9f%ac%eb%8b%9d%ec%9d%84-%eb%8f%84%ec%9e%85%ed%95%9c-%eb%b0%a9%ec%8b%9d-1/"

# Row 46: update title only
$ws.Range("D46").Value = "[Bioinformatics] 2021년 12월,  한국유전체학회 제18회 동계 워크샵"
